# Atualização automática SALDO_PECAS (14/11/2025 20:18)
# Adds one new data row (row 10) to the PRINCIPAL sheet, mirroring the
# existing rows' layout (UF, FRU, SUB1, SUB2, SUB3, DESCRICAO, MAQUINAS,
# CLIENTE, DATA_FIM, SLA, DATA_VERIFICACAO, STATUS, DATA_FIM_DT).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

$ws.Cells.Item($row, 1).Value = "DF"
$ws.Cells.Item($row, 2).Value = "DF00001"

# SUB1 / SUB2 / SUB3 are blank on this row, same as every other data row.
$ws.Cells.Item($row, 3).Value = ""
$ws.Cells.Item($row, 4).Value = ""
$ws.Cells.Item($row, 5).Value = ""

$ws.Cells.Item($row, 6).Value = "TESTE1718"
$ws.Cells.Item($row, 7).Value = "T"
$ws.Cells.Item($row, 8).Value = "T - (T 03/11/25_12H) - DF"

# DATA_FIM / SLA / DATA_VERIFICACAO are stored as plain text (not real
# dates) in this sheet, e.g. "03/11/25" - force text formatting first so
# Excel doesn't auto-convert the day/month/year-looking string into a
# date serial number, then drop back to the Normal style so no stray
# number-format style lingers on the cell (matches the other rows).
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "03/11/25"
$ws.Cells.Item($row, 9).Style = "Normal"

$ws.Cells.Item($row, 10).Value = "12H"
$ws.Cells.Item($row, 11).Value = "14/11/25"
$ws.Cells.Item($row, 12).Value = "DENTRO"

# DATA_FIM_DT is blank on this row too (only populated for a subset of
# rows in this sheet).
$ws.Cells.Item($row, 13).Value = ""
